$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N (14) - this shifts the old
# N/O/P ("Late" / "heading+Original/Principal(E4)" / "Outstanding") columns
# one to the right (to O/P/Q) and leaves the new N column blank, matching
# the "Variable Instalments" layout change used for RBI loans.
$ws.Columns.Item(14).Insert()

# Give the newly inserted column the same width as column M (13), just
# without the "best fit" flag (matches the authored column width for the
# new blank column).
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Make "Repayment schedule" the active sheet/tab (was "Acc_Periodic").
$ws.Activate()

# Move the selection on the newly active sheet to T10.
$ws.Range("T10").Select() | Out-Null
